{"js": "// Load the three tables in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst infoTable = tables.items[0];   // Date / ID No / Name / TP No table\nconst itemTable = tables.items[1];   // Rented item details table\nconst payTable = tables.items[2];    // Payment details table\n\n// ---------------------------------------------------------------------\n// Table 1: info table -- plain text replacements, cell-by-cell, in order.\n// ---------------------------------------------------------------------\ninfoTable.getCell(0, 0).value = \"\u0daf\u0dd2\u0db1\u0dba\";\ninfoTable.getCell(0, 1).value = \"2020-12-28 13:00:22\";\ninfoTable.getCell(1, 0).value = \"\u0da2\u0dcf. \u0d85\u0d82\u0d9a\u0dba\";\ninfoTable.getCell(1, 1).value = \"123\";\ninfoTable.getCell(2, 0).value = \"\u0db1\u0db8\";\ninfoTable.getCell(2, 1).value = \"Test\";\ninfoTable.getCell(3, 0).value = \"\u0daf\u0dd4. \u0d85\u0d82\u0d9a\u0dba\";\ninfoTable.getCell(3, 1).value = \"123\";\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Table 2: Rented item details -- drop the Total row and the Days /\n// Amount columns, widen the remaining 4 columns, translate headers,\n// and refresh the remaining data.\n// ---------------------------------------------------------------------\nitemTable.rows.load(\"items\");\nawait context.sync();\n\n// Remove the trailing \"Total\" row (row index 3).\nitemTable.rows.items[3].delete();\nawait context.sync();\n\n// Remove the \"Days\" column (index 2), then the \"Amount\" column which is\n// now at index 4 (it was index 5 before the Days column was removed).\nitemTable.deleteColumns(2, 1);\nawait context.sync();\nitemTable.deleteColumns(4, 1);\nawait context.sync();\n\n// Widen all 4 remaining columns from 1440 to 2160 twips (108pt) by\n// setting it on every cell of the header row (column width is shared).\nitemTable.rows.load(\"items\");\nawait context.sync();\nconst itemHeaderRow = itemTable.rows.items[0];\nitemHeaderRow.cells.load(\"items\");\nawait context.sync();\nfor (const cell of itemHeaderRow.cells.items) {\n  cell.columnWidth = 108; // 108pt * 20 = 2160 twips\n}\nawait context.sync();\n\n// Update header texts (Name, Rented date, Qty, Rate -> Sinhala labels).\nitemTable.getCell(0, 0).value = \"\u0db7\u0dcf\u0dab\u0dca\u0da9 \u0dc0\u0dbb\u0dca\u0d9c\u0dba\";\nitemTable.getCell(0, 1).value = \"\u0d9a\u0dd4\u0dbd\u0dd2\u0dba\u0da7 \u0d9c\u0dad\u0dca \u0daf\u0dd2\u0db1\u0dba\";\nitemTable.getCell(0, 2).value = \"\u0db4\u0dca\u200d\u0dbb\u0db8\u0dcf\u0dab\u0dba\";\nitemTable.getCell(0, 3).value = \"\u0daf\u0dd2\u0db1\u0d9a\u0da7 \u0d9a\u0dd4\u0dbd\u0dd2\u0dba\";\n\n// Row 1 (Concrete Mixer): update the rented date; qty/rate stay the same.\nitemTable.getCell(1, 1).value = \"2020-12-28\";\n\n// Row 2 (Poker): update the rented date and reduce qty from 2 to 1.\nitemTable.getCell(2, 1).value = \"2020-12-28\";\nitemTable.getCell(2, 2).value = \"1\";\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Table 3: Payment details -- keep only the first row, relabel it, and\n// update its amount.\n// ---------------------------------------------------------------------\npayTable.rows.load(\"items\");\nawait context.sync();\nfor (let i = payTable.rows.items.length - 1; i >= 1; i--) {\n  payTable.rows.items[i].delete();\n}\nawait context.sync();\n\npayTable.getCell(0, 0).value = \"\u0d9c\u0dd9\u0dc0\u0dd3\u0db8\u0dca\";\npayTable.getCell(0, 1).value = \"Rs.   1,000.00\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# -----------------------------------------------------------------------\n# Table 1: info table -- plain text replacements, cell-by-cell, in order.\n# -----------------------------------------------------------------------\n$infoTable = $d.Tables(1)\n$infoTable.Cell(1, 1).Range.Text = \"\u0daf\u0dd2\u0db1\u0dba\"\n$infoTable.Cell(1, 2).Range.Text = \"2020-12-28 13:00:22\"\n$infoTable.Cell(2, 1).Range.Text = \"\u0da2\u0dcf. \u0d85\u0d82\u0d9a\u0dba\"\n$infoTable.Cell(2, 2).Range.Text = \"123\"\n$infoTable.Cell(3, 1).Range.Text = \"\u0db1\u0db8\"\n$infoTable.Cell(3, 2).Range.Text = \"Test\"\n$infoTable.Cell(4, 1).Range.Text = \"\u0daf\u0dd4. \u0d85\u0d82\u0d9a\u0dba\"\n$infoTable.Cell(4, 2).Range.Text = \"123\"\n\n# -----------------------------------------------------------------------\n# Table 2: Rented item details -- drop the Total row and the Days /\n# Amount columns, widen the remaining 4 columns, translate headers, and\n# refresh the remaining data.\n# -----------------------------------------------------------------------\n$itemTable = $d.Tables(2)\n\n# Remove the trailing \"Total\" row (row 4).\n$itemTable.Rows(4).Delete()\n\n# Remove the \"Days\" column (column 3), then the \"Amount\" column, which\n# is now column 5 (it was column 6 before the Days column was removed).\n$itemTable.Columns(3).Delete()\n$itemTable.Columns(5).Delete()\n\n# Widen all 4 remaining columns from 1440 to 2160 twips (108pt).\nfor ($c = 1; $c -le $itemTable.Columns.Count; $c++) {\n  $itemTable.Columns($c).Width = 108\n}\n\n# Update header texts (Name, Rented date, Qty, Rate -> Sinhala labels).\n$itemTable.Cell(1, 1).Range.Text = \"\u0db7\u0dcf\u0dab\u0dca\u0da9 \u0dc0\u0dbb\u0dca\u0d9c\u0dba\"\n$itemTable.Cell(1, 2).Range.Text = \"\u0d9a\u0dd4\u0dbd\u0dd2\u0dba\u0da7 \u0d9c\u0dad\u0dca \u0daf\u0dd2\u0db1\u0dba\"\n$itemTable.Cell(1, 3).Range.Text = \"\u0db4\u0dca\u200d\u0dbb\u0db8\u0dcf\u0dab\u0dba\"\n$itemTable.Cell(1, 4).Range.Text = \"\u0daf\u0dd2\u0db1\u0d9a\u0da7 \u0d9a\u0dd4\u0dbd\u0dd2\u0dba\"\n\n# Row 2 (Concrete Mixer): update the rented date; qty/rate stay the same.\n$itemTable.Cell(2, 2).Range.Text = \"2020-12-28\"\n\n# Row 3 (Poker): update the rented date and reduce qty from 2 to 1.\n$itemTable.Cell(3, 2).Range.Text = \"2020-12-28\"\n$itemTable.Cell(3, 3).Range.Text = \"1\"\n\n# -----------------------------------------------------------------------\n# Table 3: Payment details -- keep only the first row, relabel it, and\n# update its amount.\n# -----------------------------------------------------------------------\n$payTable = $d.Tables(3)\nfor ($i = $payTable.Rows.Count; $i -ge 2; $i--) {\n  $payTable.Rows($i).Delete()\n}\n\n$payTable.Cell(1, 1).Range.Text = \"\u0d9c\u0dd9\u0dc0\u0dd3\u0db8\u0dca\"\n$payTable.Cell(1, 2).Range.Text = \"Rs.   1,000.00\"\n"}
